# Update the "HotStock_Top20" rankings table (Sheet1, columns A:C, rows 2-21)
# to the new ranking snapshot. Only the cells whose value actually changes
# are touched; unchanged cells (e.g. row 1 header, row 3 col A/B) are left
# as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "首开股份"
$ws.Range("C2").Value = "卧龙电驱"

$ws.Range("C3").Value = "利欧股份"

$ws.Range("A4").Value = "三花智控"
$ws.Range("B4").Value = "卧龙电驱"
$ws.Range("C4").Value = "北方铜业"

$ws.Range("A5").Value = "山子高科"
$ws.Range("B5").Value = "首开股份"
$ws.Range("C5").Value = "三花智控"

$ws.Range("A6").Value = "卧龙电驱"
$ws.Range("B6").Value = "三花智控"
$ws.Range("C6").Value = "首开股份"

$ws.Range("A7").Value = "宁德时代"
$ws.Range("C7").Value = "华胜天成"

$ws.Range("A8").Value = "华胜天成"
$ws.Range("B8").Value = "露笑科技"
$ws.Range("C8").Value = "吉视传媒"

$ws.Range("A9").Value = "供销大集"
$ws.Range("B9").Value = "华胜天成"
$ws.Range("C9").Value = "万通发展"

$ws.Range("A10").Value = "利欧股份"
$ws.Range("B10").Value = "金发科技"
$ws.Range("C10").Value = "上海建工"

$ws.Range("A11").Value = "万向钱潮"
$ws.Range("B11").Value = "万向钱潮"
$ws.Range("C11").Value = "工业富联"

$ws.Range("A12").Value = "金发科技"
$ws.Range("B12").Value = "宁德时代"
$ws.Range("C12").Value = "指南针"

$ws.Range("A13").Value = "均胜电子"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "供销大集"

$ws.Range("A14").Value = "工业富联"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "金发科技"

$ws.Range("A15").Value = "指南针"
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = "青山纸业"

$ws.Range("A16").Value = "东方财富"
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = "北方稀土"

$ws.Range("A17").Value = "大洋电机"
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = "岩山科技"

$ws.Range("A18").Value = "露笑科技"
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = "山子高科"

$ws.Range("A19").Value = "吉视传媒"
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = "中国电影"

$ws.Range("A20").Value = "青山纸业"
$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = "数据港"

$ws.Range("A21").Value = "岩山科技"
$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = "恒宝股份"
